# Applies the changes described by the target diff to the active document.
#
# Most of the diff's textual hunks are *run-merges*: two or three adjacent
# <w:r> runs that already share identical run properties (rPr) get folded
# into a single run after the underlying text was touched by an edit. Word
# performs this run-recombination automatically whenever you edit text that
# sits on a run boundary, so the trick used below is to do a "no-op" find &
# replace (replacement text identical to the search text) whose match spans
# exactly the run(s) that should be merged. That nudges the engine to
# recombine the adjacent, identically-formatted runs without altering the
# visible text at all.

$d = $word.ActiveDocument

function NoOpReplace([string]$text) {
    $rng = $d.Content.Duplicate
    # wdFindContinue = 1, wdReplaceAll = 2
    $null = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# 1) "Description: " + "Unshielded twisted pair cable is ..." -> merge the
#    (previously separate) leading-space run into the description run.
#    Touching text inside the description run (without affecting the
#    preceding bold "Description" / ":" runs) is enough to trigger the
#    merge of the space-run with the text-run that follows it.
NoOpReplace("UTP consists")

# 2) ": Coaxial cables contain a " + "center" + " conductor and a metal
#    shield ..." -> merge all three runs into one.
NoOpReplace("contain a center conductor")

# 3) "In a coaxial cable, ... Can carry high-frequency signals" + "." ->
#    merge the trailing "." run into the sentence run.
NoOpReplace("high-frequency signals.")

# 4) Picture border line color: srgbClr 000000 -> schemeClr tx1 (same
#    rendered black, just expressed as a theme color reference instead of
#    an explicit RGB value) on each of the 6 inline pictures.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    try {
        $shp.Line.ForeColor.ObjectThemeColor = 13  # wdThemeColorText1 / msoThemeColorText1
    } catch {
    }
}

# 5) DefaultParagraphFont style becomes semiHidden (hidden from the Quick
#    Style gallery until used).
$dpf = $d.Styles("Default Paragraph Font")
try {
    $dpf.SemiHidden = $true
} catch {
    try {
        $dpf.Hidden = $true
    } catch {
    }
}
